# Insert a new data row at row 423 (pushing the existing rows 423-500 down
# to 424-501) and populate it with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 423; this shifts rows 423:500
# down to 424:501 and extends the used range to A1:R501, matching the
# rest of the existing rows (same columns/format).
$ws.Rows.Item(423).Insert()

# Populate the newly-inserted row 423 with the new record.
$ws.Range("A423").Value = 5
$ws.Range("B423").Value = "Macroferia Regional de Talca"
$ws.Range("C423").Value = "Maule"
$ws.Range("D423").Value = 45180
$ws.Range("E423").Value = 7
$ws.Range("F423").Value = 100112008
$ws.Range("G423").Value = "Coliflor"
$ws.Range("H423").Value = "Sin especificar"
$ws.Range("I423").Value = "Primera"
$ws.Range("J423").Value = 3000
$ws.Range("K423").Value = 900
$ws.Range("L423").Value = 900
$ws.Range("M423").Value = 900
$ws.Range("N423").Value = "`$/unidad"
$ws.Range("O423").Value = "Región del Maule"
$ws.Range("P423").Value = 900
$ws.Range("Q423").Value = 1
$ws.Range("R423").Value = "Hortaliza"
